$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# heading: op-test -> 中文
$ws.Range("E2").Value = "中文"

# subheading: who is there -> 英文
$ws.Range("F2").Value = "英文"

# extraheading: 4444 -> 123
$ws.Range("G2").Value = 123

# thumbnail_bg_image_path: (empty) -> path
$ws.Range("I2").Value = "/Users/wenke/github/tiktoka-studio-uploader-app/tests/videos/horizon/1920x1080/bg2.jpg"

# thumbnail_local_path: [] -> (empty)
$ws.Range("K2").Value = ""

# release_date: 2023-10-13T00:00:00+00:00 -> (empty)
$ws.Range("L2").Value = ""

# release_date_hour: 19:45 -> 10:15
$ws.Range("M2").Value = "10:15"
